# Update row 5 with the new simulated values and remove row 6 entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 5 (columns A..AB)
$ws.Range("A5").Value = 474
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0.8
$ws.Range("D5").Value = 14
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 39
$ws.Range("H5").Value = 12
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 11.2
$ws.Range("L5").Value = 13.1
$ws.Range("M5").Value = 6
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 4.800000000000001
$ws.Range("Q5").Value = 1.1
$ws.Range("R5").Value = 146.72
$ws.Range("S5").Value = 5.280000000000001
$ws.Range("T5").Value = 13
$ws.Range("U5").Value = 14.2
$ws.Range("V5").Value = "Normal"
$ws.Range("W5").Value = "MTR02"
$ws.Range("X5").Value = "15x15"
$ws.Range("Y5").Value = 474
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = 25
$ws.Range("AB5").Value = "arriba"

# Remove row 6 entirely, shrinking the used range back to A1:AB5
$ws.Rows.Item(6).Delete()
